{"js": "// 1. Update the activation date from 2016 to 2023.\nconst dateResults = context.document.body.search(\"Ativa\u00e7\u00e3o: 01/01/2016\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", Word.InsertLocation.replace);\n} else {\n  throw new Error(\"Could not find the 'Ativa\u00e7\u00e3o: 01/01/2016' text to update.\");\n}\n\n// 2. Insert italicized English translations after three Portuguese paragraphs.\nconst translations = [\n  {\n    search:\n      \"Desenvolver conhecimentos de forma a tornar o aluno capaz de interpretar corretamente o desenho t\u00e9cnico, conhecer as metodologias e ferramentas utilizadas na ind\u00fastria, dando subs\u00eddios para que possa executar, interagir e modificar desenhos e projetos ao longo de sua vida profissional.\",\n    english:\n      \"Develop knowledge in order to make the student capable of correctly interpreting the technical drawing, knowing the methodologies and tools used in the industry, giving subsidies so that they can execute, interact and modify drawings and projects throughout their professional life.\",\n  },\n  {\n    search:\n      \"Contexto do desenho t\u00e9cnico na ind\u00fastria, principais ferramentas e t\u00e9cnicas utilizadas em desenhos para elabora\u00e7\u00e3o de projetos. Introdu\u00e7\u00e3o ao desenho assistido por computador (CAD).\",\n    english:\n      \"Context of the technical drawing in the industry, main tools and techniques used in drawings for the elaboration of projects. Introduction to computer-aided design (CAD).\",\n  },\n  {\n    search:\n      \"Normas do desenho t\u00e9cnico. Terminologia t\u00e9cnica e materiais para desenho. Representa\u00e7\u00e3o em perspectiva. Projeto ortogonal. Dimensionamento e escala. Corte e sec\u00e7\u00e3o. Vista Auxiliar e detalhes. Toler\u00e2ncias geom\u00e9tricas. Representa\u00e7\u00e3o de elementos de m\u00e1quinas. Utiliza\u00e7\u00e3o de software para desenho t\u00e9cnico. Desenho assistido por computador em tr\u00eas dimens\u00f5es (Modelagem de S\u00f3lidos). Desenho assistido por computador em duas dimens\u00f5es.\",\n    english:\n      \"Technical drawing standards. Technical terminology and materials for drawing. Perspective representation. Orthogonal design. Scaling and scaling. Cut and section. Auxiliary view and details. Geometric tolerances. Representation of machine elements. Use of software for technical design. Computer-aided design in three dimensions (Solid Modeling). Computer-aided design in two dimensions.\",\n  },\n];\n\nfor (const item of translations) {\n  const results = context.document.body.search(item.search, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find paragraph text starting with: ${item.search.substring(0, 40)}...`);\n  }\n\n  const paragraph = results.items[0].paragraphs.getFirst();\n  const newParagraph = paragraph.insertParagraph(item.english, Word.InsertLocation.after);\n  newParagraph.font.set({ italic: true });\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the activation date from 2016 to 2023.\n$dateRange = $d.Content\n$found = $dateRange.Find.Execute(\"Ativa\u00e7\u00e3o: 01/01/2016\")\nif ($found) {\n    $dateRange.Text = \"Ativa\u00e7\u00e3o: 01/01/2023\"\n} else {\n    throw \"Could not find the 'Ativa\u00e7\u00e3o: 01/01/2016' text to update.\"\n}\n\n# 2. Insert italicized English translations after three Portuguese paragraphs.\n$translations = @(\n    @{\n        Search  = \"Desenvolver conhecimentos de forma a tornar o aluno capaz de interpretar corretamente o desenho t\u00e9cnico, conhecer as metodologias e ferramentas utilizadas na ind\u00fastria, dando subs\u00eddios para que possa executar, interagir e modificar desenhos e projetos ao longo de sua vida profissional.\"\n        English = \"Develop knowledge in order to make the student capable of correctly interpreting the technical drawing, knowing the methodologies and tools used in the industry, giving subsidies so that they can execute, interact and modify drawings and projects throughout their professional life.\"\n    },\n    @{\n        Search  = \"Contexto do desenho t\u00e9cnico na ind\u00fastria, principais ferramentas e t\u00e9cnicas utilizadas em desenhos para elabora\u00e7\u00e3o de projetos. Introdu\u00e7\u00e3o ao desenho assistido por computador (CAD).\"\n        English = \"Context of the technical drawing in the industry, main tools and techniques used in drawings for the elaboration of projects. Introduction to computer-aided design (CAD).\"\n    },\n    @{\n        Search  = \"Normas do desenho t\u00e9cnico. Terminologia t\u00e9cnica e materiais para desenho. Representa\u00e7\u00e3o em perspectiva. Projeto ortogonal. Dimensionamento e escala. Corte e sec\u00e7\u00e3o. Vista Auxiliar e detalhes. Toler\u00e2ncias geom\u00e9tricas. Representa\u00e7\u00e3o de elementos de m\u00e1quinas. Utiliza\u00e7\u00e3o de software para desenho t\u00e9cnico. Desenho assistido por computador em tr\u00eas dimens\u00f5es (Modelagem de S\u00f3lidos). Desenho assistido por computador em duas dimens\u00f5es.\"\n        English = \"Technical drawing standards. Technical terminology and materials for drawing. Perspective representation. Orthogonal design. Scaling and scaling. Cut and section. Auxiliary view and details. Geometric tolerances. Representation of machine elements. Use of software for technical design. Computer-aided design in three dimensions (Solid Modeling). Computer-aided design in two dimensions.\"\n    }\n)\n\nforeach ($item in $translations) {\n    $matchedParagraph = $null\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.TrimEnd() -eq $item.Search) {\n            $matchedParagraph = $p\n            break\n        }\n    }\n\n    if ($matchedParagraph -eq $null) {\n        throw \"Could not find paragraph with text: $($item.Search)\"\n    }\n\n    $matchedParagraph.Range.InsertParagraphAfter()\n    $newParagraph = $matchedParagraph.Next()\n    $newRange = $newParagraph.Range\n    $startPos = $newRange.Start\n    $newRange.InsertAfter($item.English)\n    $insertedRange = $d.Range($startPos, $startPos + $item.English.Length)\n    $insertedRange.Font.Italic = 1\n}\n"}
